$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Ativação:" date 01/01/2012 -> 01/01/2023 ---
# Affects every cell that showed the old date: B8/C8 (current) and B15/C15
# (which mirrored the same value in the original file). Using a text formula
# first avoids Excel's automatic "looks like a date" literal coercion, then
# Paste-Special-Values bakes it back down to a plain cached value while
# keeping the cell's existing style untouched.
$ws.Range("B8").Formula = '="01/01/2023"'
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)

$ws.Range("C8").Formula = '="01/01/2023"'
$ws.Range("C8").Copy()
$ws.Range("C8").PasteSpecial(-4163)

$ws.Range("B15").Formula = '="01/01/2023"'
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)

$ws.Range("C15").Formula = '="01/01/2023"'
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# --- Row 11: "Objectives:" -> fill in the English objectives text ---
$ws.Range("B11").Value = "The course aims to provide students with basic knowledge of electrochemistry, both from the point of view of ionic and electrodic electrochemistry, and to present the main applications of electrochemistry"
$ws.Range("B13").Copy()
$ws.Range("B11").PasteSpecial(-4122)

$ws.Range("C11").Value = "The course aims to provide students with basic knowledge of electrochemistry, both from the point of view of ionic and electrodic electrochemistry, and to present the main applications of electrochemistry"
$ws.Range("C13").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# --- Row 14: "Short syllabus:" -> fill in the English short syllabus text ---
$ws.Range("B14").Value = "Principles of ionic electrochemistry and electrodic electrochemistry. Applications."
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("C14").Value = "Principles of ionic electrochemistry and electrodic electrochemistry. Applications."
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# --- Row 16: "Syllabus:" -> fill in the English full syllabus text ---
$ws.Range("B16").Value = "Principles of ionic electrochemistry: ionic interactions, ionic equilibrium and electrolyte conduction. Principles of electrochemical electrochemistry: interfacial phenomena, electrode potentials and electrochemical cells. Electrode processes. Electrochemical methods of chemical analysis. Applications of electrochemistry: electrochemical sources of energy, electrometallurgical processes and electroplating."
$ws.Range("B13").Copy()
$ws.Range("B16").PasteSpecial(-4122)

$ws.Range("C16").Value = "Principles of ionic electrochemistry: ionic interactions, ionic equilibrium and electrolyte conduction. Principles of electrochemical electrochemistry: interfacial phenomena, electrode potentials and electrochemical cells. Electrode processes. Electrochemical methods of chemical analysis. Applications of electrochemistry: electrochemical sources of energy, electrometallurgical processes and electroplating."
$ws.Range("C13").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$excel.CutCopyMode = $false
